# April 22, high-level features of OOP
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (the "2" row) now documents OOP instead of higher-order functions.
$ws.Range("B3").Value = "learned high-level OOP"
$ws.Range("C3").Value = "Including __slots__, @property, multi-inheritation, customized class, Enum class, metaclass;`nUse of exception handling: try, except and finally, raise Exception"

# Row 2 (the "1" row) now just recalls OOP, with its Result/Actions cell cleared.
$ws.Range("B2").Value = "Recalled OOP"
$ws.Range("C2").Value = ""

# Row 3 grew taller to fit the new multi-line reflection text.
$ws.Rows.Item(3).RowHeight = 40.5

# Selection moved to the header/body block of the first table (A1:C3).
$ws.Range("A1:C3").Select() | Out-Null

# Remember the workbook window position from this editing session.
$win = $wb.Windows.Item(1)
$win.Left = 5393
$win.Top = 1920
